$d = $word.ActiveDocument

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p w:rsidR="007A1EB8" w:rsidRDefault="007A1EB8" w:rsidP="007A1EB8"><w:pPr><w:pStyle w:val="Kop1"/></w:pPr><w:r><w:t>Hoe wordt het tafeltennis balletje gedetecteerd?</w:t></w:r></w:p><w:p w:rsidR="009C16EB" w:rsidRDefault="009C16EB" w:rsidP="009C16EB"><w:r><w:t xml:space="preserve">Om het balletje terug te kunnen slaan moet het systeem weten waar het balletje zich bevindt. Hiervoor moet het systeem uit een informatiebron (Camera, Infrarood, sonar, </w:t></w:r><w:r w:rsidR="00ED5AA1"><w:t>etc.</w:t></w:r><w:r><w:t xml:space="preserve">) de nodige informatie kunnen halen en verwerken tot informatie die door het systeem toegepast kan worden. </w:t></w:r><w:r w:rsidR="000B410B"><w:t xml:space="preserve">Dit heet Object Tracking. </w:t></w:r><w:r><w:t xml:space="preserve">Dit proces kan samengevat worden met het volgende diagram(Gebaseerd op </w:t></w:r><w:r w:rsidRPr="009C16EB"><w:t xml:space="preserve">A Survey on Object </w:t></w:r><w:r><w:t xml:space="preserve">Detection and Tracking Methods </w:t></w:r><w:sdt><w:sdtPr><w:id w:val="1621947413"/><w:citation/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> CITATION Him14 \l 1043 </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>(Himani, Darshak, &amp; Udesang, 2014)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve">): </w:t></w:r></w:p><w:p w:rsidR="009C16EB" w:rsidRDefault="009C16EB" w:rsidP="009C16EB"><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="nl-NL"/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="5486400" cy="3200400"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="2" name="Diagram 2"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/diagram"><dgm:relIds xmlns:dgm="http://schemas.openxmlformats.org/drawingml/2006/diagram" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" r:dm="rId6" r:lo="rId7" r:qs="rId8" r:cs="rId9"/></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p><w:p w:rsidR="0065071A" w:rsidRDefault="00ED5AA1" w:rsidP="0065071A"><w:pPr><w:pStyle w:val="Kop2"/></w:pPr><w:r><w:t>Nodige eigenschappen voor de beeldherkenning</w:t></w:r></w:p><w:p w:rsidR="00556EC3" w:rsidRDefault="00556EC3" w:rsidP="00556EC3"><w:r><w:t xml:space="preserve">Om het systeem te correct te laten </w:t></w:r><w:r w:rsidR="007C1433"><w:t>tafeltennissen</w:t></w:r><w:r><w:t xml:space="preserve"> zijn een aantal aspecten van de beeldherkenning van toepassing:</w:t></w:r></w:p><w:p w:rsidR="00556EC3" w:rsidRDefault="00556EC3" w:rsidP="00556EC3"><w:pPr><w:pStyle w:val="Kop3"/></w:pPr><w:r><w:t>Uitvoer tijd</w:t></w:r></w:p><w:p w:rsidR="004E13A8" w:rsidRPr="004E13A8" w:rsidRDefault="004E13A8" w:rsidP="004E13A8"><w:r><w:t>Doordat het systeem maar een beperkte tijd heeft om te reageren zal het balletje snel getraceerd moeten worden. Wanneer de beelden niet snel genoeg verwerkt worden zal de robotarm niet op tijd kunnen reageren.</w:t></w:r></w:p><w:p w:rsidR="00556EC3" w:rsidRDefault="00556EC3" w:rsidP="00556EC3"><w:pPr><w:pStyle w:val="Kop3"/></w:pPr><w:r><w:t>Accuraatheid</w:t></w:r></w:p><w:p w:rsidR="00476666" w:rsidRPr="00476666" w:rsidRDefault="00476666" w:rsidP="00476666"><w:r><w:t>Het systeem moet weten wa</w:t></w:r><w:r w:rsidR="007C1433"><w:t>ar het balletje zich bevindt en zal hiermee moeten uitrekenen waar het balletje zich zal bevinden. Echter hoeft dit niet volledig accuraat uitgevoerd te worden. Dit komt doordat de arm zal reageren met een batje waardoor een verschil van enkele centimeters weinig verschil zal maken.</w:t></w:r></w:p><w:p w:rsidR="00556EC3" w:rsidRDefault="00556EC3" w:rsidP="00556EC3"><w:pPr><w:pStyle w:val="Kop3"/></w:pPr><w:r><w:t>Resistentie tegen verschillen</w:t></w:r></w:p><w:p w:rsidR="005813FE" w:rsidRPr="005813FE" w:rsidRDefault="005813FE" w:rsidP="005813FE"><w:r><w:t>Tussen de beelden in zullen een aantal verschillen ontstaan. Deze kunnen veroorzaakt worden door beweging op de achtergrond, verschil in lichtsterkte of het vallen van schaduw. Hierdoor zal het mogelijk zijn dat het balletje niet correct wordt gevonden of dat het systeem het balletje op een andere locatie verwacht(false positive).</w:t></w:r></w:p><w:p w:rsidR="007D548D" w:rsidRDefault="007D548D" w:rsidP="007D548D"><w:pPr><w:pStyle w:val="Kop3"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Herstelmogelijkheid wanneer het balletje buiten beeld valt</w:t></w:r></w:p><w:p w:rsidR="00556EC3" w:rsidRDefault="005813FE" w:rsidP="003A08F1"><w:r><w:t xml:space="preserve">Wanneer een speler(of de arm) een punt scoort </w:t></w:r><w:r w:rsidR="00431643"><w:t>zal het balletje van de tafel vallen en waarschijnlijk buiten het beeld terecht komen. Ook kan een van de spelers zich zo positioneren dat zijn batje(of arm) het beeld van de camera blokkeert. In deze situaties moet het systeem het balletje snel terug kunnen vinden om te zorgen dat door gespeeld kan worden.</w:t></w:r></w:p><w:p w:rsidR="00DB1A8E" w:rsidRDefault="00DB1A8E" w:rsidP="00DB1A8E"><w:pPr><w:pStyle w:val="Kop2"/></w:pPr><w:r><w:t>Beschikbare beeldherkenning technieken</w:t></w:r><w:r w:rsidR="00BB6B7B"><w:t>/middelen</w:t></w:r></w:p><w:p w:rsidR="00320FD8" w:rsidRPr="00320FD8" w:rsidRDefault="000D5049" w:rsidP="00320FD8"><w:r><w:t>Er zijn een aantal hulpmiddelen en algoritmen beschikbaar om te helpen met het vaststellen van de positie van de bal. Hieronder worden deze per categorie weergegeven:</w:t></w:r></w:p><w:p w:rsidR="00DB1A8E" w:rsidRDefault="00C02E03" w:rsidP="00C02E03"><w:pPr><w:pStyle w:val="Kop3"/></w:pPr><w:r><w:t>Opvangen beeld materiaal(Invoer beeldmateriaal</w:t></w:r><w:r w:rsidR="00ED1748"><w:t>)</w:t></w:r></w:p><w:p w:rsidR="007D2E85" w:rsidRDefault="007D2E85" w:rsidP="007D2E85"><w:r><w:t xml:space="preserve">De eerste stap in het implementeren van Object Tracking is een invoerbron waaruit </w:t></w:r><w:r w:rsidR="00B35EEE"><w:t xml:space="preserve">informatie word toegediend. </w:t></w:r><w:r w:rsidR="00ED1748"><w:t>Dit zijn ruwe beelden die meerdere malen per seconde worden opgenomen. Hiervoor bestaan de volgende technieken:</w:t></w:r></w:p><w:p w:rsidR="00ED1748" w:rsidRDefault="00ED1748" w:rsidP="00ED1748"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r><w:t>Sonar/Echo</w:t></w:r></w:p><w:p w:rsidR="0039171F" w:rsidRPr="0039171F" w:rsidRDefault="0039171F" w:rsidP="0039171F"><w:r><w:t>Een manier om het balletje te detecteren is met behulp van Sonar of Echo. Dit houdt in dat een sensor een geluidsignaal verzend en de tijd meet totdat hij dit signaal terug ontvangt. Hiermee kan een afstand worden uitgerekend tot de sensor. Echter kan dit signaal worden verstoord door andere geluidsbronnen in de omgeving</w:t></w:r><w:r w:rsidR="000E5883"><w:t>, zoals het contact van het balletje op een batje, het stuiteren van het balletje of het spreken van mensen in de omgeving.</w:t></w:r></w:p><w:p w:rsidR="00ED1748" w:rsidRDefault="00ED1748" w:rsidP="00ED1748"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r><w:t>Infrarood</w:t></w:r></w:p><w:p w:rsidR="005829A8" w:rsidRPr="005829A8" w:rsidRDefault="005829A8" w:rsidP="005829A8"><w:r><w:t>Een infrarood sensor kijkt naar de warmte die een object uitstraalt. Echter zal het balletje een soortgelijke temperatuur hebben als de omgeving waar deze zich in bevindt waardoor het balletje slecht zichtbaar zal zijn op een infrarood beeld.</w:t></w:r></w:p><w:p w:rsidR="00ED1748" w:rsidRDefault="00ED1748" w:rsidP="00ED1748"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r><w:t>Camera(kleur of grijswaarden)</w:t></w:r></w:p><w:p w:rsidR="0080298D" w:rsidRPr="00FC0C4F" w:rsidRDefault="00FC0C4F" w:rsidP="00FC0C4F"><w:r><w:t xml:space="preserve">Ten slotte kan worden gewerkt met beelden uit een camera. </w:t></w:r><w:r w:rsidR="0080298D"><w:t>Iedere camera zal de nodige beelden kunnen produceren, echter zullen camera’s met meer beelden per seconden een accuratere herkenning opleveren doordat het verschil tussen de beelden minder groot is.</w:t></w:r></w:p><w:p w:rsidR="00C02E03" w:rsidRDefault="00C02E03" w:rsidP="00C02E03"><w:pPr><w:pStyle w:val="Kop3"/></w:pPr><w:r><w:t>Detecteren van de mogelijke object positie(Object Detectie</w:t></w:r><w:r w:rsidR="00D63BAB"><w:t>)</w:t></w:r></w:p><w:p w:rsidR="0050286F" w:rsidRDefault="00D63BAB" w:rsidP="0050286F"><w:r><w:t>Vervolgens moet een techniek worden toegepast om een mogelijke positie van het balletje te verkrijgen. Het is mogelijk dat de technieken meerdere mogelijkheden retourneren, dit zal in een volgende stap</w:t></w:r><w:r w:rsidR="00601B6E"><w:t xml:space="preserve"> worden opgelost.</w:t></w:r></w:p><w:p w:rsidR="008619F4" w:rsidRDefault="008619F4" w:rsidP="008619F4"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r><w:t xml:space="preserve">Frame </w:t></w:r><w:r w:rsidRPr="00B352A6"><w:t>Differencing</w:t></w:r></w:p><w:p w:rsidR="008619F4" w:rsidRDefault="008619F4" w:rsidP="008619F4"><w:r><w:t>Met Frame Differencing wordt het huidige beeld vergeleken met een voorgaand beeld. Ieder verschil dat ontstaat is een mogelijke positie van een bewegend object. Doordat dit algoritme zeer eenvoudig is, is dit ook een vrij snelle methode om mogelijke locaties vast te stellen.</w:t></w:r><w:r w:rsidR="00AC76ED"><w:t xml:space="preserve"> Echter kunnen verschillen ontstaan door veranderingen in het licht en mogelijke bewegingen op een achtergrond.</w:t></w:r></w:p><w:p w:rsidR="00BE6D68" w:rsidRDefault="00BE6D68" w:rsidP="00BE6D68"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r><w:t>Optical Flow</w:t></w:r></w:p><w:p w:rsidR="00BE6D68" w:rsidRDefault="00BE6D68" w:rsidP="00BE6D68"><w:r><w:t xml:space="preserve">Optical flow processing is een algoritme dat ogenschijnlijke beweging van een object. </w:t></w:r><w:r w:rsidR="00B15633"><w:t>Door iedere wijziging tussen beelden te registreren en te verwerken met een serie formules. Deze methode is in staat beweging zeer gedetailleerd op te vangen, echter is deze methode zeer intensief en niet goed in staat om in variërende omstandigheden te werken.</w:t></w:r></w:p><w:p w:rsidR="000C6B53" w:rsidRDefault="000C6B53" w:rsidP="000C6B53"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Background substraction</w:t></w:r></w:p><w:p w:rsidR="000C6B53" w:rsidRDefault="000C6B53" w:rsidP="000C6B53"><w:r><w:t>Dit model</w:t></w:r><w:r w:rsidR="005E3133"><w:t xml:space="preserve"> verwerkt beelden aan de hand van een vastgestelde achtergrond(template). Vervo</w:t></w:r><w:r w:rsidR="0013467D"><w:t>lgens wordt in iedere frame deze achtergrond uit het beeld gehaald.</w:t></w:r><w:r w:rsidR="005F2164"><w:t xml:space="preserve"> Hierdoor blijven alleen de bewegende objecten over in het beeld. Echter is deze methode zwak tegen verschillen in de omgeving zoals licht en beweging op de achtergrond. Hierdoor zal deze methode een stabiele achtergrond nodig hebben met een vastgestelde lichtinbreng.</w:t></w:r><w:r w:rsidR="00B979EE"><w:t xml:space="preserve"> Ook moet de template gekalibreerd worden aan de positie van de camera.</w:t></w:r></w:p><w:p w:rsidR="000C6B53" w:rsidRDefault="000C6B53" w:rsidP="000C6B53"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r><w:t>Kleur herkenning</w:t></w:r></w:p><w:p w:rsidR="00C90316" w:rsidRDefault="0037011D" w:rsidP="0050286F"><w:r><w:t>Ten slotte kan de locatie worden ingeschat op basis van een kleur. Dit houdt in dat ieder deel van het beeld, dat niet aan de kleurvereisten voldoet, genegeerd zal worden. De resterende informatie zal een mogelijke positie zijn van een object. Deze methode vereist ook dat er weinig objecten met een soortgelijke kleur aanwezig zijn in de omgeving. Ook kan een verschil in het lichtniveau een probleem opleveren voor het systeem.</w:t></w:r></w:p><w:p w:rsidR="00C90316" w:rsidRPr="00FB646F" w:rsidRDefault="00C90316" w:rsidP="00C90316"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r w:rsidRPr="00FB646F"><w:t xml:space="preserve">Canny </w:t></w:r><w:r w:rsidR="00FB646F" w:rsidRPr="00FB646F"><w:t xml:space="preserve">Edge Detection </w:t></w:r></w:p><w:p w:rsidR="00C90316" w:rsidRPr="00FB646F" w:rsidRDefault="00FB646F" w:rsidP="00C90316"><w:r w:rsidRPr="00FB646F"><w:t xml:space="preserve">Objecten als een geheel herkennen is een ingewikkeld process voor een computer. </w:t></w:r><w:r><w:t>Een manier om het eenvoudiger te maken is het toepassen van Canny Edge Detection</w:t></w:r><w:sdt><w:sdtPr><w:id w:val="-495659628"/><w:citation/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> CITATION Can86 \l 1043 </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> (Canny, 1986)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r w:rsidR="00606721"><w:t>. Dit algoritme gebruikt een Gaussian Filter(Blur) om eventuele verstoringen en minieme veranderingen in het beeld te verwijderen. Vervolgens wordt de intensiteit van iedere pixel berekend(grijswaarde) en deze worden vergeleken met nabije pixels. Waar de intensiteit sterk verschilt van de aangrenzende pixels kan gesproken worden over een rand. Deze gegevens worden in het zwart-wit getekend op een afbeelding van dezelfde grootte als het origineel. Hierd</w:t></w:r><w:r w:rsidR="00D01D51"><w:t>oor raken de kleuren uit het beeld verloren.</w:t></w:r><w:r w:rsidR="000E6574"><w:t xml:space="preserve"> Daarentegen is het vrij ongevoelig voor verschillen in lichtsterkte.</w:t></w:r></w:p><w:p w:rsidR="00C02E03" w:rsidRDefault="00C02E03" w:rsidP="00C02E03"><w:pPr><w:pStyle w:val="Kop3"/></w:pPr><w:r><w:t>Herkenning van het object(Object Classificatie)</w:t></w:r><w:r w:rsidR="009929F7"><w:t xml:space="preserve"> en vaststellen positie</w:t></w:r></w:p><w:p w:rsidR="00571284" w:rsidRDefault="00571284" w:rsidP="00571284"><w:r><w:t xml:space="preserve">Rafeal Nieto heeft een overzicht gemaakt van mogelijke Object </w:t></w:r><w:r w:rsidR="00FD3A2D"><w:t>classificatie</w:t></w:r><w:r><w:t xml:space="preserve"> technieken die worden gebruikt in zijn Master Thesis</w:t></w:r><w:sdt><w:sdtPr><w:id w:val="-278648226"/><w:citation/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> CITATION Raf13 \l 1043 </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> (Nieto, 2013)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve"> heeft omschreven.</w:t></w:r><w:r w:rsidR="00500B7A"><w:t xml:space="preserve"> Hieronder staat een beknopt overzicht van deze methoden.</w:t></w:r></w:p><w:p w:rsidR="00500B7A" w:rsidRDefault="00500B7A" w:rsidP="00500B7A"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r><w:t>Template Matching</w:t></w:r></w:p><w:p w:rsidR="0073340D" w:rsidRDefault="0073340D" w:rsidP="00B352A6"><w:r><w:t>Bij Template Matching wordt op een beeld het juiste object opgespoord door deze te vergelijken met een vooraf gesteld beeld(template)</w:t></w:r><w:r w:rsidR="00B352A6"><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="00C050A6"><w:t xml:space="preserve">Dit wordt gedaan door een Convolutie(Verschil tussen waarden </w:t></w:r><w:sdt><w:sdtPr><w:id w:val="413052405"/><w:citation/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r w:rsidR="00C050A6"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="00C050A6"><w:instrText xml:space="preserve"> CITATION Con15 \l 1043 </w:instrText></w:r><w:r w:rsidR="00C050A6"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00C050A6"><w:rPr><w:noProof/></w:rPr><w:t>(Convolution, 2015)</w:t></w:r><w:r w:rsidR="00C050A6"><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r w:rsidR="00C050A6"><w:t>) uit te rekenen en de locatie met de hoogste convolutie waarde is het object dat gevonden dient te worden.</w:t></w:r></w:p><w:p w:rsidR="00C050A6" w:rsidRDefault="00C050A6" w:rsidP="00B352A6"><w:r><w:t>Doordat het algoritme uit weinig stappen bestaat is het eenvoud</w:t></w:r><w:r w:rsidR="00D5047C"><w:t>ig toe te passen. Deze methode verwerkt zijn gegevens snel genoeg om in een real-time applicatie te kunnen draaien.</w:t></w:r></w:p><w:p w:rsidR="00D5047C" w:rsidRDefault="00D5047C" w:rsidP="00B352A6"><w:r><w:t>Echter kan het algoritme niet goed tegen transformaties van het doelobject (Vervormingen, rotaties en verandering in formaat). Deze problemen zullen niet veel voorkomen bij de ping pong bal doordat deze ten alle tijden rond zullen zijn.</w:t></w:r></w:p><w:p w:rsidR="00D5047C" w:rsidRPr="0073340D" w:rsidRDefault="00D5047C" w:rsidP="00B352A6"><w:r><w:t>Verder kunnen problemen ontstaan wanneer de kleur van het object afwijkt van de template. Dit kan komen door bijvoorbeeld de belichting van het object. Dit probleem kan echter opgelost worden in de Object Detectie stap.</w:t></w:r></w:p><w:p w:rsidR="00500B7A" w:rsidRDefault="00500B7A" w:rsidP="00500B7A"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Mean-Shift</w:t></w:r></w:p><w:p w:rsidR="00C358AF" w:rsidRDefault="00D8266C" w:rsidP="00C358AF"><w:r><w:t>Mean Shift beschrijft een proces waarbij de nieuwe locatie wordt uitgerekend aan de hand van een eerdere positie in combinatie met een herkenningspunt. Dit herkenningspunt kan een template, een kleur combinatie of een andere herkenningspunt</w:t></w:r><w:r w:rsidR="005E5296"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>zijn.</w:t></w:r></w:p><w:p w:rsidR="00A2722D" w:rsidRDefault="005E5296" w:rsidP="00C358AF"><w:r><w:t xml:space="preserve">Het algoritme rekent locaties uit die overeenkomen met het herkenningspunt. Echter is deze herkenning niet zo strikt als bij template matches, waarbij iedere de volledige overeenkomst telt, maar in plaats hiervan zullen alle overeenkomsten gemarkeerd worden. Vervolgens wordt de </w:t></w:r><w:r w:rsidRPr="005E5296"><w:t xml:space="preserve">Epanechnikov </w:t></w:r><w:r><w:t xml:space="preserve">Kernel methode </w:t></w:r><w:sdt><w:sdtPr><w:id w:val="805280149"/><w:citation/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> CITATION Str15 \l 1043 </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>(Struijker)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r w:rsidR="00A2722D"><w:t xml:space="preserve"> toegepast om de overeenkomsten dicht bij de eerdere locatie prioriteit te geven over de anderen overeenkomsten.</w:t></w:r><w:r><w:t xml:space="preserve"> Vervolgens wordt het algoritme iteratief uitgevoerd over het resultaat totdat de locaties van de objecten samenvallen.</w:t></w:r></w:p><w:p w:rsidR="00A2722D" w:rsidRDefault="00A2722D" w:rsidP="00C358AF"><w:r><w:t>Deze methode kan ook doorberekend worden wanneer het object buiten beeld valt(Obstructie) door met de eerder uitgerekende snelheid en locatie de verwachte locatie van het object aan te passen. Hierdoor kan ieder frame een voorspelling gemaakt worden totdat een nieuwe cluster overeenkomsten in de buurt van de verwachte locatie komt waardoor het object teruggevonden kan worden.</w:t></w:r></w:p><w:p w:rsidR="00A2722D" w:rsidRDefault="00A2722D" w:rsidP="00C358AF"><w:r><w:t>Het algoritme presteert goed wanneer een specifiek kenmerk het object omschrijft</w:t></w:r><w:r w:rsidR="00D90B6E"><w:t>. Hierdoor kan het object eenvoudig herkent worden.</w:t></w:r></w:p><w:p w:rsidR="005E5296" w:rsidRPr="00C358AF" w:rsidRDefault="00A2722D" w:rsidP="00C358AF"><w:r><w:t xml:space="preserve">Dit algoritme presteert echter </w:t></w:r><w:r w:rsidR="00D90B6E"><w:t>minder goed wanneer het object buiten het beeldbereik valt. Dit komt doordat het algoritme zal aannemen dat het object door beweegt met als gevolg dat de voorspelde locatie ver buiten het beeld zal vallen.</w:t></w:r><w:r w:rsidR="005E5296"><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w:rsidR="00500B7A" w:rsidRDefault="00500B7A" w:rsidP="00500B7A"><w:pPr><w:pStyle w:val="Kop4"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00231970"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Particle Filter-b</w:t></w:r><w:r w:rsidRPr="00B352A6"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>ased Color Tracking (PFC)</w:t></w:r></w:p><w:p w:rsidR="00C61164" w:rsidRDefault="00C61164" w:rsidP="00C61164"><w:r w:rsidRPr="00C61164"><w:t>Dit algoritme werkt op een verglijkbare manier als de bovenstaande Mean-Shift methode</w:t></w:r><w:r w:rsidR="009E610B"><w:t xml:space="preserve">, echter werkt deze alleen met kleur. Vervolgens wordt hier ook de bovenstaande </w:t></w:r><w:r w:rsidR="009E610B" w:rsidRPr="005E5296"><w:t xml:space="preserve">Epanechnikov </w:t></w:r><w:r w:rsidR="009E610B"><w:t xml:space="preserve">Kernel methode </w:t></w:r><w:sdt><w:sdtPr><w:id w:val="-899592356"/><w:citation/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r w:rsidR="009E610B"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="009E610B"><w:instrText xml:space="preserve"> CITATION Str15 \l 1043 </w:instrText></w:r><w:r w:rsidR="009E610B"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="009E610B"><w:rPr><w:noProof/></w:rPr><w:t>(Struijker)</w:t></w:r><w:r w:rsidR="009E610B"><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r w:rsidR="009E610B"><w:t xml:space="preserve"> toegepast om de overeenkomsten dicht bij de vorige positie een hogere waarde te geven. </w:t></w:r></w:p><w:p w:rsidR="009E610B" w:rsidRDefault="009E610B" w:rsidP="00C61164"><w:r><w:t>Echter gebruikt deze methode het gemiddelde van alle potentiele locaties om de daadwerkelijke locatie te bepalen. Deze vergelijking van locaties wordt herhaaldelijk toegepast to</w:t></w:r><w:r w:rsidR="005B7E4D"><w:t>t</w:t></w:r><w:r><w:t>dat de locaties samenvallen op een locatie</w:t></w:r><w:r w:rsidR="005B7E4D"><w:t>.</w:t></w:r><w:r w:rsidR="005A1DCC"><w:t xml:space="preserve"> Deze uiteindelijke locatie is het resultaat van het algoritme.</w:t></w:r></w:p><w:p w:rsidR="005A1DCC" w:rsidRPr="00C61164" w:rsidRDefault="005A1DCC" w:rsidP="00C61164"><w:r><w:t>Dit algoritme presteert beter als anderen in complexe situaties doordat iedere mogelijke locati</w:t></w:r><w:r w:rsidR="00580087"><w:t>e mee</w:t></w:r><w:r><w:t>genomen wordt in het eind resultaat.</w:t></w:r><w:r w:rsidR="00580087"><w:t xml:space="preserve"> Hierdoor is dit een van de meest gebruikte algoritmes voor beeldherkenning.</w:t></w:r></w:p><w:p w:rsidR="00500B7A" w:rsidRPr="00231970" w:rsidRDefault="00500B7A" w:rsidP="00500B7A"><w:pPr><w:pStyle w:val="Kop4"/></w:pPr><w:r w:rsidRPr="00231970"><w:t>Lucas-Kanade Tracking</w:t></w:r></w:p><w:p w:rsidR="00ED33CD" w:rsidRDefault="00ED33CD" w:rsidP="00ED33CD"><w:r w:rsidRPr="00231970"><w:t xml:space="preserve">Lucas-Kanade Tracking </w:t></w:r><w:r w:rsidR="00231970" w:rsidRPr="00231970"><w:t>is een vorm van Optical Flow herkenning</w:t></w:r><w:sdt><w:sdtPr><w:id w:val="612946974"/><w:citation/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r w:rsidR="00231970"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="00231970"><w:instrText xml:space="preserve"> CITATION Roj15 \l 1043 </w:instrText></w:r><w:r w:rsidR="00231970"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00231970"><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> (Rojas)</w:t></w:r><w:r w:rsidR="00231970"><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r w:rsidR="00231970" w:rsidRPr="00231970"><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="00231970"><w:t>Het werkt door de verschillen in grijswaarden(Intensiteit) te meten. Hierbij worden beelden vergleken en probeert het algoritme in te schatten welke richting het object op beweegt zodat de verschillen in intensiteit verklaard kunnen worden.</w:t></w:r></w:p><w:p w:rsidR="00231970" w:rsidRPr="00231970" w:rsidRDefault="00231970" w:rsidP="00ED33CD"><w:r><w:t>Echter gaat het algoritme er van uit dat het verschil in tijd en afgelegde afstand tussen de beelden niet groot is.</w:t></w:r><w:r w:rsidR="003321F7"><w:t xml:space="preserve"> Dit wil zeggen dat het algoritme is gebouwd is voor langzaam bewegende objecten.</w:t></w:r><w:r w:rsidR="00CC3486"><w:t xml:space="preserve"> Hierdoor is deze niet geschikt voor het bijhouden van tafeltennis.</w:t></w:r></w:p><w:p w:rsidR="00500B7A" w:rsidRDefault="00500B7A" w:rsidP="00500B7A"><w:pPr><w:pStyle w:val="Kop4"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00B352A6"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t>Incremental Learning for Robust Visual Tracking</w:t></w:r></w:p><w:p w:rsidR="00C90316" w:rsidRPr="003C3D06" w:rsidRDefault="003C3D06" w:rsidP="00C90316"><w:r w:rsidRPr="003C3D06"><w:t>Incremental Learning gebruikt een aantal templates om een inschatting te maken van de mogelijke transformaties van een object.</w:t></w:r><w:r><w:t xml:space="preserve"> Vervolgens gaat het algoritme van ieder nieuw beeld het object registreren in een nieuwe template. Door het beeld te vergelijken met alle templates zal het algoritme eventuele veranderingen kunnen opvangen doordat hij het nieuwe object zal kennen. Daarentegen kan het algoritme niet goed omgaan met het verlies van het object en zal het algoritme steeds intensiever worden doordat de hoeveelheid templates groter wordt.</w:t></w:r></w:p><w:p w:rsidR="00500B7A" w:rsidRPr="00B352A6" w:rsidRDefault="00500B7A" w:rsidP="00500B7A"><w:pPr><w:pStyle w:val="Kop4"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00B352A6"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Tracking Learning Detection</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">De Tracking Learning Detection is een tracker die zich aanpast aan het beeld dat  deze ontvangt. </w:t></w:r><w:r><w:t xml:space="preserve">Dit wordt gerealiseerd door informatie van de frames bij te houden en </w:t></w:r><w:r><w:t xml:space="preserve">te gebruiken. Deze tracker gaat er hierdoor vanuit dat de beweging tussen de frames </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>klein is(lage snelheid) en dat het object in beeld is. Wanneer het object buiten beeld valt zal de tracker het object niet meer terug kunnen vinden.</w:t></w:r></w:p><w:p w:rsidR="00500B7A" w:rsidRPr="00B352A6" w:rsidRDefault="00500B7A" w:rsidP="00500B7A"><w:pPr><w:pStyle w:val="Kop4"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00B352A6"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Corrected Background-Weighted Histogram Tracker</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Het doel van de Background weighted histogram tracker is om de Mean-Shift methode te verbeteren door </w:t></w:r><w:r><w:t>de invloeden vanuit de achtergrond te verminderen.</w:t></w:r><w:sdt><w:sdtPr><w:id w:val="-1489621609"/><w:citation/></w:sdtPr><w:sdtContent><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> CITATION Yan13 \l 1043 </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> (Yang, Jia, Rong, Zhu, Wang, &amp; Yue, 2013)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve"> Daarentegen is hij door deze methode minder efficiënt wanne</w:t></w:r><w:r><w:t>er er kleurverschillen optreden. Om dit probleem op te lossen kan een Kalman Filter toegepast worden.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Deze methode kan de hoeveelheid iteraties van de Mean-Shift methode verkleinen. Hierdoor zal deze variant het object sneller kunnen vinden. Ook zal het een hogere </w:t></w:r><w:r><w:t>precisie</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>kunnen behalen doordat minder objecten meetellen in de berekening.</w:t></w:r></w:p><w:p><w:r><w:t>Het algoritme presteert goed wanneer er een duidelijk verschil is tussen het doel en de achtergrond. Echter kan het algoritme problemen krijgen wanneer er meerdere soortgelijke objecten bij elkaar aanwezig zijn.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="00500B7A" w:rsidRPr="00B352A6" w:rsidRDefault="00500B7A" w:rsidP="00500B7A"><w:pPr><w:pStyle w:val="Kop4"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00B352A6"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Scale and Orientation Adaptive Mean-Shift Tracking</w:t></w:r></w:p><w:p w:rsidR="00C02E03" w:rsidRPr="00B352A6" w:rsidRDefault="00C02E03" w:rsidP="00C02E03"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p w:rsidR="0065071A" w:rsidRPr="00B352A6" w:rsidRDefault="0065071A" w:rsidP="0065071A"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p w:rsidR="007A1EB8" w:rsidRPr="00B352A6" w:rsidRDefault="007A1EB8" w:rsidP="007A1EB8"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:sectPr w:rsidR="007A1EB8" w:rsidRPr="00B352A6"><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1417" w:right="1417" w:bottom="1417" w:left="1417" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml)
